$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: Hide existing rows that transition from visible to hidden
$hideRows = @(1591,1604,1623,1651,1653,1670,1694,1718,1742,1763,1777,1789,1813,1837,1853,1859,1880,1895,1930,1947,1972,1982,2008,2029,2053,2071,2087,2110,2123,2138,2197,2206,2223,2236,2256,2270,2297,2313,2363,2394,2419,2437,2446,2461,2476,2560,2561,2562,2563,2564,2565,2566,2567,2568,2569,2570,2571,2572,2573,2574,2575,2576,2577,2578,2579,2580,2581,2582,2583,2584,2585,2586,2587,2588,2589,2590,2591,2592,2593,2594,2595,2596,2597,2598,2599,2600,2601,2602,2603,2604,2605,2606,2607,2608,2609,2610,2611,2612,2613,2614,2615,2616,2617,2618,2619,2620,2621,2622,2623,2624,2625,2626,2627,2628,2629,2630,2631,2632,2633,2634,2635,2636,2637,2638,2639,2640,2641,2642,2643,2644,2645,2646,2647,2648,2649,2650,2651,2652,2653,2654,2655,2656,2657,2658,2659,2660,2661,2662)
foreach ($r in $hideRows) {
    $ws.Rows.Item($r).Hidden = $true
}

# Step 2: Insert new rows 2663-2709 using row 2662 as a formatting template
$lastRow = 2662
for ($i = 0; $i -lt 47; $i++) {
    $ws.Rows.Item($lastRow).Copy()
    $ws.Rows.Item($lastRow + 1).Insert(-4121)
}

# Step 3: Populate new rows with values, then set row height + hidden status
$ws.Cells.Item(2663, 1).Value = 46013.63638091435
$ws.Cells.Item(2663, 2).Value = "Lunes"
$ws.Cells.Item(2663, 3).Value = "Dagoberto Espinoza"
$ws.Cells.Item(2663, 4).Value = "Normal"
$ws.Cells.Item(2663, 5).Value = "Normal"
$ws.Cells.Item(2663, 6).Value = "6 a 8"
$ws.Cells.Item(2663, 7).Value = "No me duele nada"
$ws.Cells.Item(2663, 8).Value = "Nada"
$ws.Rows.Item(2663).RowHeight = 22.5
$ws.Rows.Item(2663).Hidden = $true

$ws.Cells.Item(2664, 1).Value = 46013.64290060185
$ws.Cells.Item(2664, 2).Value = "Lunes"
$ws.Cells.Item(2664, 3).Value = "Rodrigo Aguirre"
$ws.Cells.Item(2664, 4).Value = "Normal"
$ws.Cells.Item(2664, 5).Value = "Normal"
$ws.Cells.Item(2664, 6).Value = "6 a 8"
$ws.Cells.Item(2664, 7).Value = "Adolorido de una zona"
$ws.Cells.Item(2664, 8).Value = "10 isquiotibial derecho"
$ws.Rows.Item(2664).RowHeight = 22.5
$ws.Rows.Item(2664).Hidden = $true

$ws.Cells.Item(2665, 1).Value = 46013.644082893516
$ws.Cells.Item(2665, 2).Value = "Lunes"
$ws.Cells.Item(2665, 3).Value = "Jonathan Dos Santos"
$ws.Cells.Item(2665, 4).Value = "Normal"
$ws.Cells.Item(2665, 5).Value = "Normal"
$ws.Cells.Item(2665, 6).Value = "6 a 8"
$ws.Cells.Item(2665, 7).Value = "Adolorido de una zona"
$ws.Cells.Item(2665, 8).Value = "9 Isquiotibial izquierdo"
$ws.Rows.Item(2665).RowHeight = 22.5
$ws.Rows.Item(2665).Hidden = $true

$ws.Cells.Item(2666, 1).Value = 46013.646062245374
$ws.Cells.Item(2666, 2).Value = "Lunes"
$ws.Cells.Item(2666, 3).Value = "Santiago Naveda"
$ws.Cells.Item(2666, 4).Value = "Normal"
$ws.Cells.Item(2666, 5).Value = "Normal"
$ws.Cells.Item(2666, 6).Value = "6 a 8"
$ws.Cells.Item(2666, 7).Value = "Adolorido de una zona"
$ws.Cells.Item(2666, 8).Value = "10 isquiotibial derecho"
$ws.Rows.Item(2666).RowHeight = 22.5
$ws.Rows.Item(2666).Hidden = $true

$ws.Cells.Item(2667, 1).Value = 46013.64692101852
$ws.Cells.Item(2667, 2).Value = "Lunes"
$ws.Cells.Item(2667, 3).Value = "Luis Ángel Malagón"
$ws.Cells.Item(2667, 4).Value = "Normal"
$ws.Cells.Item(2667, 5).Value = "Normal"
$ws.Cells.Item(2667, 6).Value = "6 a 8"
$ws.Cells.Item(2667, 7).Value = "Adolorido de una zona"
$ws.Cells.Item(2667, 8).Value = "10 isquiotibial derecho"
$ws.Rows.Item(2667).RowHeight = 22.5
$ws.Rows.Item(2667).Hidden = $true

$ws.Cells.Item(2668, 1).Value = 46013.65180277778
$ws.Cells.Item(2668, 2).Value = "Lunes"
$ws.Cells.Item(2668, 3).Value = "Alan Cervantes"
$ws.Cells.Item(2668, 4).Value = "Normal"
$ws.Cells.Item(2668, 5).Value = "Normal"
$ws.Cells.Item(2668, 6).Value = "6 a 8"
$ws.Cells.Item(2668, 7).Value = "Adolorido de una zona"
$ws.Cells.Item(2668, 8).Value = "4 rodilla izquierda, 7 tobillo derecho"
$ws.Rows.Item(2668).RowHeight = 22.5
$ws.Rows.Item(2668).Hidden = $true

$ws.Cells.Item(2669, 1).Value = 46013.65762892361
$ws.Cells.Item(2669, 2).Value = "Lunes"
$ws.Cells.Item(2669, 3).Value = "Henry Martín"
$ws.Cells.Item(2669, 4).Value = "Normal"
$ws.Cells.Item(2669, 5).Value = "Normal"
$ws.Cells.Item(2669, 6).Value = "Más de 8"
$ws.Cells.Item(2669, 7).Value = "Adolorido de una zona"
$ws.Cells.Item(2669, 8).Value = "18 aductor derecho"
$ws.Rows.Item(2669).RowHeight = 22.5
$ws.Rows.Item(2669).Hidden = $true

$ws.Cells.Item(2670, 1).Value = 46013.6579621875
$ws.Cells.Item(2670, 2).Value = "Lunes"
$ws.Cells.Item(2670, 3).Value = "Álvaro Fidalgo"
$ws.Cells.Item(2670, 4).Value = "Normal"
$ws.Cells.Item(2670, 5).Value = "Normal"
$ws.Cells.Item(2670, 6).Value = "6 a 8"
$ws.Cells.Item(2670, 7).Value = "Adolorido de una zona"
$ws.Cells.Item(2670, 8).Value = "15 espalda baja"
$ws.Rows.Item(2670).RowHeight = 22.5
$ws.Rows.Item(2670).Hidden = $true

$ws.Cells.Item(2671, 1).Value = 46013.66350510417
$ws.Cells.Item(2671, 2).Value = "Lunes"
$ws.Cells.Item(2671, 3).Value = "Erick Sánchez"
$ws.Cells.Item(2671, 4).Value = "Peor que lo normal"
$ws.Cells.Item(2671, 5).Value = "Normal"
$ws.Cells.Item(2671, 6).Value = "6 a 8"
$ws.Cells.Item(2671, 7).Value = "No me duele nada"
$ws.Cells.Item(2671, 8).Value = "Nada"
$ws.Rows.Item(2671).RowHeight = 22.5
$ws.Rows.Item(2671).Hidden = $true

$ws.Cells.Item(2672, 1).Value = 46013.663840381945
$ws.Cells.Item(2672, 2).Value = "Lunes"
$ws.Cells.Item(2672, 3).Value = "Miguel Vázquez"
$ws.Cells.Item(2672, 4).Value = "Normal"
$ws.Cells.Item(2672, 5).Value = "Normal"
$ws.Cells.Item(2672, 6).Value = "6 a 8"
$ws.Cells.Item(2672, 7).Value = "No me duele nada"
$ws.Cells.Item(2672, 8).Value = "Nada"
$ws.Rows.Item(2672).RowHeight = 22.5
$ws.Rows.Item(2672).Hidden = $true

$ws.Cells.Item(2673, 1).Value = 46013.66417862268
$ws.Cells.Item(2673, 2).Value = "Lunes"
$ws.Cells.Item(2673, 3).Value = "Víctor Dávila"
$ws.Cells.Item(2673, 4).Value = "Normal"
$ws.Cells.Item(2673, 5).Value = "Normal"
$ws.Cells.Item(2673, 6).Value = "6 a 8"
$ws.Cells.Item(2673, 7).Value = "No me duele nada"
$ws.Cells.Item(2673, 8).Value = "Nada"
$ws.Rows.Item(2673).RowHeight = 22.5
$ws.Rows.Item(2673).Hidden = $true

$ws.Cells.Item(2674, 1).Value = 46013.664499108796
$ws.Cells.Item(2674, 2).Value = "Lunes"
$ws.Cells.Item(2674, 3).Value = "Alexis Gutiérrez"
$ws.Cells.Item(2674, 4).Value = "Normal"
$ws.Cells.Item(2674, 5).Value = "Normal"
$ws.Cells.Item(2674, 6).Value = "6 a 8"
$ws.Cells.Item(2674, 7).Value = "No me duele nada"
$ws.Cells.Item(2674, 8).Value = "Nada"
$ws.Rows.Item(2674).RowHeight = 22.5
$ws.Rows.Item(2674).Hidden = $true

$ws.Cells.Item(2675, 1).Value = 46013.6647893287
$ws.Cells.Item(2675, 2).Value = "Lunes"
$ws.Cells.Item(2675, 3).Value = "Isaías Violante"
$ws.Cells.Item(2675, 4).Value = "Normal"
$ws.Cells.Item(2675, 5).Value = "Normal"
$ws.Cells.Item(2675, 6).Value = "6 a 8"
$ws.Cells.Item(2675, 7).Value = "No me duele nada"
$ws.Cells.Item(2675, 8).Value = "Nada"
$ws.Rows.Item(2675).RowHeight = 22.5
$ws.Rows.Item(2675).Hidden = $true

$ws.Cells.Item(2676, 1).Value = 46013.665156006944
$ws.Cells.Item(2676, 2).Value = "Lunes"
$ws.Cells.Item(2676, 3).Value = "José Raúl Zúñiga"
$ws.Cells.Item(2676, 4).Value = "Peor que lo normal"
$ws.Cells.Item(2676, 5).Value = "Normal"
$ws.Cells.Item(2676, 6).Value = "6 a 8"
$ws.Cells.Item(2676, 7).Value = "Normal"
$ws.Cells.Item(2676, 8).Value = "Nada"
$ws.Rows.Item(2676).RowHeight = 22.5
$ws.Rows.Item(2676).Hidden = $true

$ws.Cells.Item(2677, 1).Value = 46013.66688496528
$ws.Cells.Item(2677, 2).Value = "Lunes"
$ws.Cells.Item(2677, 3).Value = "Sebastián Cáceres"
$ws.Cells.Item(2677, 4).Value = "Normal"
$ws.Cells.Item(2677, 5).Value = "Normal"
$ws.Cells.Item(2677, 6).Value = "6 a 8"
$ws.Cells.Item(2677, 7).Value = "No me duele nada"
$ws.Cells.Item(2677, 8).Value = "Nada"
$ws.Rows.Item(2677).RowHeight = 22.5
$ws.Rows.Item(2677).Hidden = $true

$ws.Cells.Item(2678, 1).Value = 46013.66868671296
$ws.Cells.Item(2678, 2).Value = "Lunes"
$ws.Cells.Item(2678, 3).Value = "Israel Reyes"
$ws.Cells.Item(2678, 4).Value = "Normal"
$ws.Cells.Item(2678, 5).Value = "Normal"
$ws.Cells.Item(2678, 6).Value = "6 a 8"
$ws.Cells.Item(2678, 7).Value = "Adolorido de una zona"
$ws.Cells.Item(2678, 8).Value = "5 espinilla derecha, 6 espinilla izquierda"
$ws.Rows.Item(2678).RowHeight = 22.5
$ws.Rows.Item(2678).Hidden = $true

$ws.Cells.Item(2679, 1).Value = 46013.670060358796
$ws.Cells.Item(2679, 2).Value = "Lunes"
$ws.Cells.Item(2679, 3).Value = "Ramón Juárez"
$ws.Cells.Item(2679, 4).Value = "Normal"
$ws.Cells.Item(2679, 5).Value = "Normal"
$ws.Cells.Item(2679, 6).Value = "6 a 8"
$ws.Cells.Item(2679, 7).Value = "Adolorido de una zona"
$ws.Cells.Item(2679, 8).Value = "13 pantorrilla izquierda, 14 pantorrilla derecha"
$ws.Rows.Item(2679).RowHeight = 22.5
$ws.Rows.Item(2679).Hidden = $true

$ws.Cells.Item(2680, 1).Value = 46013.67331460648
$ws.Cells.Item(2680, 2).Value = "Lunes"
$ws.Cells.Item(2680, 3).Value = "Igor Lichnovsky"
$ws.Cells.Item(2680, 4).Value = "Normal"
$ws.Cells.Item(2680, 5).Value = "Normal"
$ws.Cells.Item(2680, 6).Value = "6 a 8"
$ws.Cells.Item(2680, 7).Value = "No me duele nada"
$ws.Cells.Item(2680, 8).Value = "Nada"
$ws.Rows.Item(2680).RowHeight = 22.5
$ws.Rows.Item(2680).Hidden = $true

$ws.Cells.Item(2681, 1).Value = 46013.6736204051
$ws.Cells.Item(2681, 2).Value = "Lunes"
$ws.Cells.Item(2681, 3).Value = "Brian Rodríguez"
$ws.Cells.Item(2681, 4).Value = "Normal"
$ws.Cells.Item(2681, 5).Value = "Normal"
$ws.Cells.Item(2681, 6).Value = "6 a 8"
$ws.Cells.Item(2681, 7).Value = "No me duele nada"
$ws.Cells.Item(2681, 8).Value = "Nada"
$ws.Rows.Item(2681).RowHeight = 22.5
$ws.Rows.Item(2681).Hidden = $true

$ws.Cells.Item(2682, 1).Value = 46013.67488200231
$ws.Cells.Item(2682, 2).Value = "Lunes"
$ws.Cells.Item(2682, 3).Value = "Alejandro Zendejas"
$ws.Cells.Item(2682, 4).Value = "Normal"
$ws.Cells.Item(2682, 5).Value = "Normal"
$ws.Cells.Item(2682, 6).Value = "6 a 8"
$ws.Cells.Item(2682, 7).Value = "Adolorido de una zona"
$ws.Cells.Item(2682, 8).Value = "9 Isquiotibial izquierdo"
$ws.Rows.Item(2682).RowHeight = 22.5
$ws.Rows.Item(2682).Hidden = $true

$ws.Cells.Item(2683, 1).Value = 46013.68261574074
$ws.Cells.Item(2683, 2).Value = "Lunes"
$ws.Cells.Item(2683, 3).Value = "Rodolfo Cota"
$ws.Cells.Item(2683, 4).Value = "Normal"
$ws.Cells.Item(2683, 5).Value = "Normal"
$ws.Cells.Item(2683, 6).Value = "6 a 8"
$ws.Cells.Item(2683, 7).Value = "No me duele nada"
$ws.Cells.Item(2683, 8).Value = "Nada"
$ws.Rows.Item(2683).RowHeight = 22.5
$ws.Rows.Item(2683).Hidden = $true

$ws.Cells.Item(2684, 1).Value = 46013.68291307871
$ws.Cells.Item(2684, 2).Value = "Lunes"
$ws.Cells.Item(2684, 3).Value = "Néstor Araujo"
$ws.Cells.Item(2684, 4).Value = "Normal"
$ws.Cells.Item(2684, 5).Value = "Normal"
$ws.Cells.Item(2684, 6).Value = "6 a 8"
$ws.Cells.Item(2684, 7).Value = "No me duele nada"
$ws.Cells.Item(2684, 8).Value = "Nada"
$ws.Rows.Item(2684).RowHeight = 22.5
$ws.Rows.Item(2684).Hidden = $true

$ws.Cells.Item(2685, 1).Value = 46013.68362329861
$ws.Cells.Item(2685, 2).Value = "Lunes"
$ws.Cells.Item(2685, 3).Value = "Cristian Borja"
$ws.Cells.Item(2685, 4).Value = "Normal"
$ws.Cells.Item(2685, 5).Value = "Me desperté mucho"
$ws.Cells.Item(2685, 6).Value = "Menos de 6"
$ws.Cells.Item(2685, 7).Value = "No me duele nada"
$ws.Cells.Item(2685, 8).Value = "Nada"
$ws.Rows.Item(2685).RowHeight = 22.5
$ws.Rows.Item(2685).Hidden = $true

$ws.Cells.Item(2686, 1).Value = 46013.68426207176
$ws.Cells.Item(2686, 2).Value = "Lunes"
$ws.Cells.Item(2686, 3).Value = "Kevin Álvarez"
$ws.Cells.Item(2686, 4).Value = "Normal"
$ws.Cells.Item(2686, 5).Value = "Normal"
$ws.Cells.Item(2686, 6).Value = "6 a 8"
$ws.Cells.Item(2686, 7).Value = "No me duele nada"
$ws.Cells.Item(2686, 8).Value = "Nada"
$ws.Rows.Item(2686).RowHeight = 22.5
$ws.Rows.Item(2686).Hidden = $true

$ws.Cells.Item(2687, 1).Value = 46013.68832883102
$ws.Cells.Item(2687, 2).Value = "Lunes"
$ws.Cells.Item(2687, 3).Value = "Allan Saint-Maximin"
$ws.Cells.Item(2687, 4).Value = "Normal"
$ws.Cells.Item(2687, 5).Value = "Normal"
$ws.Cells.Item(2687, 6).Value = "6 a 8"
$ws.Cells.Item(2687, 7).Value = "No me duele nada"
$ws.Cells.Item(2687, 8).Value = "Nada"
$ws.Rows.Item(2687).RowHeight = 22.5
$ws.Rows.Item(2687).Hidden = $true

$ws.Cells.Item(2688, 1).Value = 46013.688577152774
$ws.Cells.Item(2688, 2).Value = "Lunes"
$ws.Cells.Item(2688, 3).Value = "Ralph Orquin"
$ws.Cells.Item(2688, 4).Value = "Normal"
$ws.Cells.Item(2688, 5).Value = "Normal"
$ws.Cells.Item(2688, 6).Value = "6 a 8"
$ws.Cells.Item(2688, 7).Value = "No me duele nada"
$ws.Cells.Item(2688, 8).Value = "Nada"
$ws.Rows.Item(2688).RowHeight = 22.5
$ws.Rows.Item(2688).Hidden = $true

$ws.Cells.Item(2689, 1).Value = 46014.355254062495
$ws.Cells.Item(2689, 2).Value = "Martes"
$ws.Cells.Item(2689, 3).Value = "Igor Lichnovsky"
$ws.Cells.Item(2689, 4).Value = "Normal"
$ws.Cells.Item(2689, 5).Value = "Muy buena noche"
$ws.Cells.Item(2689, 6).Value = "Más de 8"
$ws.Cells.Item(2689, 7).Value = "Muy adolorido en general"
$ws.Cells.Item(2689, 8).Value = "Nada"
$ws.Rows.Item(2689).RowHeight = 22.5

$ws.Cells.Item(2690, 1).Value = 46014.366081238426
$ws.Cells.Item(2690, 2).Value = "Martes"
$ws.Cells.Item(2690, 3).Value = "Néstor Araujo"
$ws.Cells.Item(2690, 4).Value = "Peor que lo normal"
$ws.Cells.Item(2690, 5).Value = "Normal"
$ws.Cells.Item(2690, 6).Value = "6 a 8"
$ws.Cells.Item(2690, 7).Value = "No me duele nada"
$ws.Cells.Item(2690, 8).Value = "Nada"
$ws.Rows.Item(2690).RowHeight = 22.5

$ws.Cells.Item(2691, 1).Value = 46014.366509178246
$ws.Cells.Item(2691, 2).Value = "Martes"
$ws.Cells.Item(2691, 3).Value = "Sebastián Cáceres"
$ws.Cells.Item(2691, 4).Value = "Normal"
$ws.Cells.Item(2691, 5).Value = "Normal"
$ws.Cells.Item(2691, 6).Value = "6 a 8"
$ws.Cells.Item(2691, 7).Value = "Adolorido de una zona"
$ws.Cells.Item(2691, 8).Value = "2 muslo derecho"
$ws.Rows.Item(2691).RowHeight = 22.5

$ws.Cells.Item(2692, 1).Value = 46014.36729189815
$ws.Cells.Item(2692, 2).Value = "Martes"
$ws.Cells.Item(2692, 3).Value = "Alan Cervantes"
$ws.Cells.Item(2692, 4).Value = "Peor que lo normal"
$ws.Cells.Item(2692, 5).Value = "Normal"
$ws.Cells.Item(2692, 6).Value = "6 a 8"
$ws.Cells.Item(2692, 7).Value = "Adolorido de una zona"
$ws.Cells.Item(2692, 8).Value = "4 rodilla izquierda, 7 tobillo derecho"
$ws.Rows.Item(2692).RowHeight = 22.5

$ws.Cells.Item(2693, 1).Value = 46014.36779369213
$ws.Cells.Item(2693, 2).Value = "Martes"
$ws.Cells.Item(2693, 3).Value = "Erick Sánchez"
$ws.Cells.Item(2693, 4).Value = "Normal"
$ws.Cells.Item(2693, 5).Value = "Normal"
$ws.Cells.Item(2693, 6).Value = "6 a 8"
$ws.Cells.Item(2693, 7).Value = "Adolorido de una zona"
$ws.Cells.Item(2693, 8).Value = "9 Isquiotibial izquierdo, 10 isquiotibial derecho"
$ws.Rows.Item(2693).RowHeight = 22.5

$ws.Cells.Item(2694, 1).Value = 46014.368286620374
$ws.Cells.Item(2694, 2).Value = "Martes"
$ws.Cells.Item(2694, 3).Value = "Álvaro Fidalgo"
$ws.Cells.Item(2694, 4).Value = "Normal"
$ws.Cells.Item(2694, 5).Value = "Normal"
$ws.Cells.Item(2694, 6).Value = "6 a 8"
$ws.Cells.Item(2694, 7).Value = "No me duele nada"
$ws.Cells.Item(2694, 8).Value = "Nada"
$ws.Rows.Item(2694).RowHeight = 22.5

$ws.Cells.Item(2695, 1).Value = 46014.36884717592
$ws.Cells.Item(2695, 2).Value = "Martes"
$ws.Cells.Item(2695, 3).Value = "Israel Reyes"
$ws.Cells.Item(2695, 4).Value = "Peor que lo normal"
$ws.Cells.Item(2695, 5).Value = "Normal"
$ws.Cells.Item(2695, 6).Value = "6 a 8"
$ws.Cells.Item(2695, 7).Value = "Adolorido de una zona"
$ws.Cells.Item(2695, 8).Value = "13 pantorrilla izquierda, 14 pantorrilla derecha"
$ws.Rows.Item(2695).RowHeight = 22.5

$ws.Cells.Item(2696, 1).Value = 46014.369300185186
$ws.Cells.Item(2696, 2).Value = "Martes"
$ws.Cells.Item(2696, 3).Value = "Jonathan Dos Santos"
$ws.Cells.Item(2696, 4).Value = "Normal"
$ws.Cells.Item(2696, 5).Value = "Normal"
$ws.Cells.Item(2696, 6).Value = "6 a 8"
$ws.Cells.Item(2696, 7).Value = "Adolorido de una zona"
$ws.Cells.Item(2696, 8).Value = "9 Isquiotibial izquierdo, 15 espalda baja"
$ws.Rows.Item(2696).RowHeight = 22.5

$ws.Cells.Item(2697, 1).Value = 46014.36972478009
$ws.Cells.Item(2697, 2).Value = "Martes"
$ws.Cells.Item(2697, 3).Value = "Kevin Álvarez"
$ws.Cells.Item(2697, 4).Value = "Normal"
$ws.Cells.Item(2697, 5).Value = "Normal"
$ws.Cells.Item(2697, 6).Value = "6 a 8"
$ws.Cells.Item(2697, 7).Value = "Adolorido de una zona"
$ws.Cells.Item(2697, 8).Value = "21 Plantas de los pies o empeine"
$ws.Rows.Item(2697).RowHeight = 22.5

$ws.Cells.Item(2698, 1).Value = 46014.37022012731
$ws.Cells.Item(2698, 2).Value = "Martes"
$ws.Cells.Item(2698, 3).Value = "Luis Ángel Malagón"
$ws.Cells.Item(2698, 4).Value = "Cansado"
$ws.Cells.Item(2698, 5).Value = "Normal"
$ws.Cells.Item(2698, 6).Value = "6 a 8"
$ws.Cells.Item(2698, 7).Value = "No me duele nada"
$ws.Cells.Item(2698, 8).Value = "Nada"
$ws.Rows.Item(2698).RowHeight = 22.5

$ws.Cells.Item(2699, 1).Value = 46014.37053568287
$ws.Cells.Item(2699, 2).Value = "Martes"
$ws.Cells.Item(2699, 3).Value = "Miguel Vázquez"
$ws.Cells.Item(2699, 4).Value = "Normal"
$ws.Cells.Item(2699, 5).Value = "Normal"
$ws.Cells.Item(2699, 6).Value = "6 a 8"
$ws.Cells.Item(2699, 7).Value = "No me duele nada"
$ws.Cells.Item(2699, 8).Value = "Nada"
$ws.Rows.Item(2699).RowHeight = 22.5

$ws.Cells.Item(2700, 1).Value = 46014.37097525463
$ws.Cells.Item(2700, 2).Value = "Martes"
$ws.Cells.Item(2700, 3).Value = "Ramón Juárez"
$ws.Cells.Item(2700, 4).Value = "Mejor que lo normal"
$ws.Cells.Item(2700, 5).Value = "Normal"
$ws.Cells.Item(2700, 6).Value = "6 a 8"
$ws.Cells.Item(2700, 7).Value = "No me duele nada"
$ws.Cells.Item(2700, 8).Value = "Nada"
$ws.Rows.Item(2700).RowHeight = 22.5

$ws.Cells.Item(2701, 1).Value = 46014.37143585648
$ws.Cells.Item(2701, 2).Value = "Martes"
$ws.Cells.Item(2701, 3).Value = "Alejandro Zendejas"
$ws.Cells.Item(2701, 4).Value = "Cansado"
$ws.Cells.Item(2701, 5).Value = "Normal"
$ws.Cells.Item(2701, 6).Value = "6 a 8"
$ws.Cells.Item(2701, 7).Value = "No me duele nada"
$ws.Cells.Item(2701, 8).Value = "Nada"
$ws.Rows.Item(2701).RowHeight = 22.5

$ws.Cells.Item(2702, 1).Value = 46014.37370295139
$ws.Cells.Item(2702, 2).Value = "Martes"
$ws.Cells.Item(2702, 3).Value = "Rodrigo Aguirre"
$ws.Cells.Item(2702, 4).Value = "Normal"
$ws.Cells.Item(2702, 5).Value = "Normal"
$ws.Cells.Item(2702, 6).Value = "6 a 8"
$ws.Cells.Item(2702, 7).Value = "Adolorido de una zona"
$ws.Cells.Item(2702, 8).Value = "23 gluteo derecho"
$ws.Rows.Item(2702).RowHeight = 22.5

$ws.Cells.Item(2703, 1).Value = 46014.374090625
$ws.Cells.Item(2703, 2).Value = "Martes"
$ws.Cells.Item(2703, 3).Value = "Cristian Borja"
$ws.Cells.Item(2703, 4).Value = "Peor que lo normal"
$ws.Cells.Item(2703, 5).Value = "Normal"
$ws.Cells.Item(2703, 6).Value = "6 a 8"
$ws.Cells.Item(2703, 7).Value = "No me duele nada"
$ws.Cells.Item(2703, 8).Value = "Nada"
$ws.Rows.Item(2703).RowHeight = 22.5

$ws.Cells.Item(2704, 1).Value = 46014.374461886575
$ws.Cells.Item(2704, 2).Value = "Martes"
$ws.Cells.Item(2704, 3).Value = "Santiago Naveda"
$ws.Cells.Item(2704, 4).Value = "Normal"
$ws.Cells.Item(2704, 5).Value = "Normal"
$ws.Cells.Item(2704, 6).Value = "6 a 8"
$ws.Cells.Item(2704, 7).Value = "Adolorido de una zona"
$ws.Cells.Item(2704, 8).Value = "9 Isquiotibial izquierdo, 10 isquiotibial derecho"
$ws.Rows.Item(2704).RowHeight = 22.5

$ws.Cells.Item(2705, 1).Value = 46014.374851342596
$ws.Cells.Item(2705, 2).Value = "Martes"
$ws.Cells.Item(2705, 3).Value = "Ralph Orquin"
$ws.Cells.Item(2705, 4).Value = "Normal"
$ws.Cells.Item(2705, 5).Value = "Normal"
$ws.Cells.Item(2705, 6).Value = "6 a 8"
$ws.Cells.Item(2705, 7).Value = "No me duele nada"
$ws.Cells.Item(2705, 8).Value = "Nada"
$ws.Rows.Item(2705).RowHeight = 22.5

$ws.Cells.Item(2706, 1).Value = 46014.375244062496
$ws.Cells.Item(2706, 2).Value = "Martes"
$ws.Cells.Item(2706, 3).Value = "Alexis Gutiérrez"
$ws.Cells.Item(2706, 4).Value = "Normal"
$ws.Cells.Item(2706, 5).Value = "Normal"
$ws.Cells.Item(2706, 6).Value = "6 a 8"
$ws.Cells.Item(2706, 7).Value = "No me duele nada"
$ws.Cells.Item(2706, 8).Value = "Nada"
$ws.Rows.Item(2706).RowHeight = 22.5

$ws.Cells.Item(2707, 1).Value = 46014.37572994213
$ws.Cells.Item(2707, 2).Value = "Martes"
$ws.Cells.Item(2707, 3).Value = "Isaías Violante"
$ws.Cells.Item(2707, 4).Value = "Normal"
$ws.Cells.Item(2707, 5).Value = "Normal"
$ws.Cells.Item(2707, 6).Value = "6 a 8"
$ws.Cells.Item(2707, 7).Value = "Adolorido de una zona"
$ws.Cells.Item(2707, 8).Value = "21 Plantas de los pies o empeine"
$ws.Rows.Item(2707).RowHeight = 22.5

$ws.Cells.Item(2708, 1).Value = 46014.37598351852
$ws.Cells.Item(2708, 2).Value = "Martes"
$ws.Cells.Item(2708, 3).Value = "José Raúl Zúñiga"
$ws.Cells.Item(2708, 4).Value = "Normal"
$ws.Cells.Item(2708, 5).Value = "Normal"
$ws.Cells.Item(2708, 6).Value = "6 a 8"
$ws.Cells.Item(2708, 7).Value = "No me duele nada"
$ws.Cells.Item(2708, 8).Value = "Nada"
$ws.Rows.Item(2708).RowHeight = 22.5

$ws.Cells.Item(2709, 1).Value = 46014.377127314816
$ws.Cells.Item(2709, 2).Value = "Martes"
$ws.Cells.Item(2709, 3).Value = "Brian Rodríguez"
$ws.Cells.Item(2709, 4).Value = "Normal"
$ws.Cells.Item(2709, 5).Value = "Mejor que normal"
$ws.Cells.Item(2709, 6).Value = "Más de 8"
$ws.Cells.Item(2709, 7).Value = "Adolorido de una zona"
$ws.Cells.Item(2709, 8).Value = "14 pantorrilla derecha"
$ws.Rows.Item(2709).RowHeight = 22.5

# Step 4: Update defined name _FilterDatabase range
foreach ($dn in $wb.Names) {
    if ($dn.Name -like "*_FilterDatabase*") {
        $dn.RefersTo = "='Respuestas de formulario 1'!`$A`$1:`$H`$2709"
    }
}

# Step 5: Resize table to new range
$tbl = $ws.ListObjects.Item(1)
$newRange = $ws.Range("A1:H2709")
$tbl.Resize($newRange)

# Step 6: Update conditional formatting ranges
$rngD = $ws.Range("D1:D2762")
$fcD = $rngD.FormatConditions.Item(1)
$fcD.ModifyAppliesToRange($ws.Range("D1:D2809"))

$rngG = $ws.Range("G2:G2762")
$fcG = $rngG.FormatConditions.Item(1)
$fcG.ModifyAppliesToRange($ws.Range("G2:G2809"))
